# Update countries & provincias Spain
# Refresh Covid-19 "Pais" sheet data: update the "last updated" timestamp,
# refresh the case counters for several countries and fix the ranking
# order for a few country pairs whose totals now place them differently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 5 de Julio de 2020 a las 07:45"

# India
$ws.Cells.Item(7,4).Value2 = 409083
$ws.Cells.Item(7,5).Value2 = 245542

# Pakistan
$ws.Cells.Item(15,2).Value2 = 228474
$ws.Cells.Item(15,3).Value2 = 3191
$ws.Cells.Item(15,4).Value2 = 129830
$ws.Cells.Item(15,5).Value2 = 93932
$ws.Cells.Item(15,7).Value2 = 93
$ws.Cells.Item(15,8).Value2 = 4712

# Uzbekistan / Sudan swap places in the ranking
$ws.Cells.Item(71,1).Value2 = "Uzbekistan"
$ws.Cells.Item(71,2).Value2 = 9829
$ws.Cells.Item(71,3).Value2 = 121
$ws.Cells.Item(71,4).Value2 = 6425
$ws.Cells.Item(71,5).Value2 = 3373
$ws.Cells.Item(71,7).Value2 = 0
$ws.Cells.Item(71,8).Value2 = 31

$ws.Cells.Item(72,1).Value2 = "Sudan"
$ws.Cells.Item(72,2).Value2 = 9767
$ws.Cells.Item(72,3).Value2 = 104
$ws.Cells.Item(72,4).Value2 = 4673
$ws.Cells.Item(72,5).Value2 = 4486
$ws.Cells.Item(72,7).Value2 = 4
$ws.Cells.Item(72,8).Value2 = 608

# El Salvador
$ws.Cells.Item(75,2).Value2 = 8449
$ws.Cells.Item(75,3).Value2 = 87
$ws.Cells.Item(75,5).Value2 = 946

# Kirguistan / Senegal / Finlandia reshuffle
$ws.Cells.Item(79,1).Value2 = "Kirguistan"
$ws.Cells.Item(79,2).Value2 = 7377
$ws.Cells.Item(79,3).Value2 = 283
$ws.Cells.Item(79,4).Value2 = 2802
$ws.Cells.Item(79,5).Value2 = 4487
$ws.Cells.Item(79,7).Value2 = 10
$ws.Cells.Item(79,8).Value2 = 88

$ws.Cells.Item(80,1).Value2 = "Senegal"
$ws.Cells.Item(80,2).Value2 = 7272
$ws.Cells.Item(80,4).Value2 = 4713
$ws.Cells.Item(80,5).Value2 = 2430
$ws.Cells.Item(80,8).Value2 = 129

$ws.Cells.Item(81,1).Value2 = "Finlandia"
$ws.Cells.Item(81,2).Value2 = 7248
$ws.Cells.Item(81,4).Value2 = 6700
$ws.Cells.Item(81,5).Value2 = 219
$ws.Cells.Item(81,8).Value2 = 329

# Costa Rica
$ws.Cells.Item(94,5).Value2 = 2881
$ws.Cells.Item(94,7).Value2 = 1
$ws.Cells.Item(94,8).Value2 = 19

# Tailandia
$ws.Cells.Item(100,2).Value2 = 3190
$ws.Cells.Item(100,3).Value2 = 5
$ws.Cells.Item(100,4).Value2 = 3071

# Fiyi / Dominica swap places in the ranking
$ws.Cells.Item(205,1).Value2 = "Fiyi"
$ws.Cells.Item(206,1).Value2 = "Dominica"
